$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A22 used to hold the literal value 6; it now holds the formula =2*3
# which evaluates to the same cached value (6).
$ws.Range("A22").Formula = "=2*3"

# B22's label changes from "six" to "six, as formula" to reflect that
# its neighbour is now computed rather than a literal.
$ws.Range("B22").Value = "six, as formula"

# Reflect the cursor/selection ending up on B23, and the view scrolled
# so row 12 is at the top, as in the saved workbook.
$ws.Range("B23").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
